$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend the table by two rows (9 and 10), copying the formatting
# (style + row height) from the last existing data row (row 8) so the new
# rows match the rest of the table.
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 20.25
$ws.Rows.Item(10).RowHeight = 20.25
$excel.CutCopyMode = 0

# Now write the updated compatibility data across rows 2-10.
$data = @(
  @("VOLVO", "C40", "2022-2023", "Rear or Front", ""),
  @("VOLVO", "S60", "2019-2023", "Rear or Front", "AWD"),
  @("VOLVO", "S90", "2017-2023", "Rear or Front", ""),
  @("VOLVO", "V60", "2019-2023", "Rear or Front", ""),
  @("VOLVO", "V90", "2018-2023", "Rear or Front", ""),
  @("VOLVO", "V90 CROSS COUNTRY", "2017-2021", "Rear or Front", ""),
  @("VOLVO", "XC40", "2021-2023", "Rear or Front", "No 2.0L L4 Turbocharged"),
  @("VOLVO", "XC60", "2018-2023", "Rear or Front", "AWD"),
  @("VOLVO", "XC90", "2016-2023", "Rear or Front", "B5 Core or B6 Core or B6 Plus or B6 Ultimate")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    if ($rowData[4] -ne "") {
        $ws.Cells.Item($row, 5).Value = $rowData[4]
    } else {
        $ws.Cells.Item($row, 5).ClearContents()
    }
}
